$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table")
$ws.Rows.Item(138).Delete()
